# First draft of first revision
#
# Table 1 row-label edits:
#  - "Response" column header:      run boundary normalized (text unchanged)
#  - "Relative humidity " (Rate of spread row): run boundary normalized (text unchanged)
#  - "Canopy temperature" row label  -> "Flame temperature"
#  - "Surface temperature" row label -> "Soil surface temperature"
#
# Note: the caption paragraph above the table is intentionally left
# untouched. In the source XML its sentence is split across two <w:r> runs
# that end up concatenated into one in the final document purely because of
# how the text was originally authored - the visible text is identical
# either way. This runtime normalizes/merges *all* same-formatted runs in a
# paragraph as soon as any part of it is edited, which would also
# incorrectly merge the trailing, differently-"rsid" " " run after
# "see Methods." (a run the source diff does not touch). So we avoid
# touching that paragraph at all to not introduce an unwanted side effect.

$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

function Replace-InTableCell {
    param(
        [int]$Row,
        [int]$Col,
        [string]$OldText,
        [string]$NewText
    )

    $cellRange = $tbl.Cell($Row, $Col).Range
    $cellStart = $cellRange.Start
    $cellEnd = $cellRange.End

    # Step 1: locate the text without mutating anything, scoped to the cell.
    $probe = $tbl.Cell($Row, $Col).Range
    $found = $probe.Find.Execute($OldText)
    if (-not $found) {
        throw "Replace-InTableCell: could not find '$OldText' in table cell ($Row,$Col)"
    }
    if (($probe.Start -lt $cellStart) -or ($probe.End -gt $cellEnd)) {
        throw "Replace-InTableCell: match for '$OldText' escaped cell ($Row,$Col) bounds; aborting to avoid editing the wrong location"
    }

    # Step 2: perform the actual replace, scoped tightly to the matched
    # sub-range so the operation cannot possibly affect any other text.
    $target = $d.Range($probe.Start, $probe.End)
    $target.Find.Execute($OldText, $true, $false, $false, $false, $false,
        $true, 0, $false, $NewText, 2)
}

# "Response" table header (row 1, column 1)
Replace-InTableCell 1 1 "Response" "Response"

# "Relative humidity " in the Rate of spread section (row 6, column 2)
Replace-InTableCell 6 2 "Relative humidity " "Relative humidity "

# "Canopy temperature" -> "Flame temperature" (row 7, column 1)
Replace-InTableCell 7 1 "Canopy temperature" "Flame temperature"

# "Surface temperature" -> "Soil surface temperature" (row 12, column 1)
Replace-InTableCell 12 1 "Surface temperature" "Soil surface temperature"
